$wb = $excel.ActiveWorkbook
$wsFeature = $wb.Worksheets.Item("feature_importance")

# Add the new worksheet right after "feature_importance" and name it.
$ws = $wb.Worksheets.Add($null, $wsFeature)
$ws.Name = "cv_scores"

# Header cell: "Cross-validation R^2 scores" in B1.
$ws.Range("B1").Value = "Cross-validation R^2 scores"

# Cross-validation fold index (col A) and R^2 score (col B).
$scores = @(0.310646068017594, 0.2424582917579408, 0.4400098584980944, 0.3365234053737414, 0.5022688053484423)

for ($i = 0; $i -lt $scores.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $scores[$i]
}

# Header + fold-index cells use the bold / centered / bordered look
# (matching the header style used on the feature_importance sheet).
$headerRange = $ws.Range("B1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$indexRange = $ws.Range("A2:A6")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

$ws.Range("A1").Select()

# Keep "feature_importance" as the active/selected sheet, as before.
$wsFeature.Activate()
